{"js": "// Receipt template (\"recibo_1.docx\") value correction: R$ 550,00 -> R$ 480,00\n// (and its written-out Portuguese form \"quinhentos e cinquenta reais\" ->\n// \"quatrocentos e oitenta reais\"), per commit message\n// \"mudei 550 pra 480 recibo frutal\".\n\nconst body = context.document.body;\n\n// 1) Numeric amount: \"R$ 550,00\" -> \"R$ 480,00\"\nconst amountHits = body.search(\"R$ 550,00\", { matchCase: true });\namountHits.load(\"text\");\nawait context.sync();\n\nfor (const hit of amountHits.items) {\n  hit.insertText(\"R$ 480,00\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Written-out amount: \"quinhentos e cinquenta reais\" -> \"quatrocentos e oitenta reais\"\nconst wordsHits = body.search(\"quinhentos e cinquenta reais\", { matchCase: true });\nwordsHits.load(\"text\");\nawait context.sync();\n\nfor (const hit of wordsHits.items) {\n  hit.insertText(\"quatrocentos e oitenta reais\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Receipt template (\"recibo_1.docx\") value correction: R$ 550,00 -> R$ 480,00\n# (and its written-out Portuguese form \"quinhentos e cinquenta reais\" ->\n# \"quatrocentos e oitenta reais\"), per commit message\n# \"mudei 550 pra 480 recibo frutal\".\n\n$d = $word.ActiveDocument\n\n# 1) Numeric amount: \"R$ 550,00\" -> \"R$ 480,00\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"R$ 550,00\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"R$ 480,00\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Written-out amount: \"quinhentos e cinquenta reais\" -> \"quatrocentos e oitenta reais\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"quinhentos e cinquenta reais\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"quatrocentos e oitenta reais\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
